$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values that look numeric must be forced to stay as text (matching the
# source data, which stores prices/percentages as inline strings) by switching
# the cell to a text number format before assigning, then restoring "Normal"
# style so no stray formatting is left behind.
function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

$ws.Range('D2').Value = '30.520.10'
$ws.Range('E2').Value = '  +0.00%  '
$ws.Range('D3').Value = '1.918.26'
$ws.Range('E3').Value = '  -0.26%  '
Set-TextValue 'D4' '0.9999'
$ws.Range('E4').Value = '  +0.05%  '
Set-TextValue 'D5' '245.73'
$ws.Range('E5').Value = '  +1.03%  '
Set-TextValue 'D6' '0.9997'
$ws.Range('E6').Value = '  +0.00%  '
Set-TextValue 'D7' '0.4794'
Set-TextValue 'D8' '0.2900'
$ws.Range('E8').Value = '  +0.72%  '
Set-TextValue 'D9' '0.06719'
$ws.Range('E9').Value = '  -0.64%  '
Set-TextValue 'D10' '110.95'
$ws.Range('E10').Value = '  +4.07%  '
Set-TextValue 'D11' '19.03'
$ws.Range('E11').Value = '  +3.82%  '
$ws.Range('D12').Value = '1.910.40'
$ws.Range('E12').Value = '  -0.60%  '
Set-TextValue 'D13' '0.07567'
$ws.Range('E13').Value = '  -2.54%  '
Set-TextValue 'D14' '5.279'
$ws.Range('E14').Value = '  -0.31%  '
Set-TextValue 'D15' '0.6678'
$ws.Range('E15').Value = '  +0.88%  '
Set-TextValue 'D16' '300.23'
$ws.Range('E16').Value = '  +2.48%  '
$ws.Range('D17').Value = '30.512.79'
$ws.Range('E17').Value = '  +0.02%  '
Set-TextValue 'D18' '5.626'
$ws.Range('E18').Value = '  +6.12%  '
Set-TextValue 'D19' '13.00'
$ws.Range('E19').Value = '  +0.38%  '
Set-TextValue 'D20' '0.9998'
$ws.Range('E20').Value = '  -0.04%  '
Set-TextValue 'D21' '0.000007571'
$ws.Range('E21').Value = '  -0.45%  '
$ws.Range('D22').Value = '2.160.47'
$ws.Range('E22').Value = '  +0.07%  '
Set-TextValue 'D23' '1.000'
$ws.Range('E23').Value = '  +0.01%  '
Set-TextValue 'D24' '6.479'
$ws.Range('E24').Value = '  +4.31%  '
Set-TextValue 'D25' '9.477'
$ws.Range('E25').Value = '  +1.03%  '
$ws.Range('E26').Value = '  -2.31%  '
$ws.Range('E27').Value = '  -5.10%  '
$ws.Range('E28').Value = '  +0.17%  '
Set-TextValue 'D29' '0.1079'
$ws.Range('E29').Value = '  +0.68%  '
Set-TextValue 'D30' '1.398'
$ws.Range('E30').Value = '  +2.30%  '
Set-TextValue 'D31' '4.172'
$ws.Range('E31').Value = '  -0.21%  '
Set-TextValue 'D32' '4.053'
$ws.Range('E32').Value = '  +0.93%  '
Set-TextValue 'D33' '0.05001'
$ws.Range('E33').Value = '  -0.81%  '
Set-TextValue 'D34' '0.7382'
$ws.Range('E34').Value = '  -0.78%  '
Set-TextValue 'D35' '1.138'
$ws.Range('E35').Value = '  -1.53%  '
$ws.Range('E37').Value = '  -0.19%  '
$ws.Range('E38').Value = '  -3.71%  '
Set-TextValue 'D39' '2.684'
$ws.Range('E39').Value = '  -0.11%  '
Set-TextValue 'D40' '111.07'
$ws.Range('E40').Value = '  +0.70%  '
Set-TextValue 'D41' '2.021'
$ws.Range('E41').Value = '  -2.66%  '
Set-TextValue 'D42' '0.4439'
$ws.Range('E42').Value = '  +3.91%  '
Set-TextValue 'D43' '72.50'
$ws.Range('E43').Value = '  +7.25%  '
Set-TextValue 'D44' '0.8642'
$ws.Range('E44').Value = '  -1.25%  '
Set-TextValue 'D45' '5.873'
$ws.Range('E45').Value = '  -0.71%  '
Set-TextValue 'D46' '0.9997'
$ws.Range('E46').Value = '  -0.04%  '
$ws.Range('B47').Value = 'Aptos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue 'D47' '7.283'
$ws.Range('E47').Value = '  +1.01%  '
$ws.Range('B48').Value = 'BitcoinSV'
$ws.Range('C48').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
Set-TextValue 'D48' '49.37'
$ws.Range('E48').Value = '  -0.57%  '
Set-TextValue 'D49' '9.315'
$ws.Range('E49').Value = '  -0.51%  '
Set-TextValue 'D50' '0.2563'
$ws.Range('E50').Value = '  +3.67%  '
$ws.Range('E51').Value = '  +0.81%  '
